$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.402088999999998
$ws.Range("H2").Value = 28.206267
$ws.Range("I2").Value = 0.4769398566373552
$ws.Range("J2").Value = 0.4769398566373552
$ws.Range("M2").Value = 1.332201
$ws.Range("N2").Value = 3.996603
$ws.Range("O2").Value = 0.4977000102551081
$ws.Range("P2").Value = 0.5003087024806611
$ws.Range("Q2").Value = 12.525472367889
$ws.Range("R2").Value = 112.729251311001
$ws.Range("S2").Value = 0.2373729715394815
$ws.Range("T2").Value = 0.2386171608355477
$ws.Range("G3").Value = 9.402088999999998
$ws.Range("H3").Value = 28.206267
$ws.Range("I3").Value = 0.4769398566373552
$ws.Range("J3").Value = 0.4769398566373552
$ws.Range("O3").Value = 0.3400184143635273
$ws.Range("P3").Value = 0.341800619257677
$ws.Range("Q3").Value = 8.557145199777665
$ws.Range("R3").Value = 77.01430679799898
$ws.Range("S3").Value = 0.1621683338006015
$ws.Range("T3").Value = 0.1630183383473157
$ws.Range("G4").Value = 9.402088999999998
$ws.Range("H4").Value = 28.206267
$ws.Range("I4").Value = 0.4769398566373552
$ws.Range("J4").Value = 0.4769398566373552
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07237833333333334
$ws.Range("N4").Value = 0.217135
$ws.Range("O4").Value = 0.02703998664034003
$ws.Range("P4").Value = 0.02718171660110807
$ws.Range("Q4").Value = 0.6805075316716666
$ws.Range("R4").Value = 6.124567785045
$ws.Range("S4").Value = 0.01289644735171977
$ws.Range("T4").Value = 0.0129640440188897
$ws.Range("G5").Value = 9.402088999999998
$ws.Range("H5").Value = 28.206267
$ws.Range("I5").Value = 0.4769398566373552
$ws.Range("J5").Value = 0.4769398566373552
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.0418705
$ws.Range("N5").Value = 0.083741
$ws.Range("O5").Value = 0.01564249559892726
$ws.Range("P5").Value = 0.01048299044324218
$ws.Range("Q5").Value = 0.3936701674744999
$ws.Range("R5").Value = 2.362021004847
$ws.Range("S5").Value = 0.007460529608402828
$ws.Range("T5").Value = 0.004999755959130691
$ws.Range("G6").Value = 9.402088999999998
$ws.Range("H6").Value = 28.206267
$ws.Range("I6").Value = 0.4769398566373552
$ws.Range("J6").Value = 0.4769398566373552
$ws.Range("M6").Value = 0.3201326666666667
$ws.Range("N6").Value = 0.960398
$ws.Range("O6").Value = 0.1195990931420972
$ws.Range("P6").Value = 0.1202259712173118
$ws.Range("Q6").Value = 3.009915823807333
$ws.Range("R6").Value = 27.089242414266
$ws.Range("S6").Value = 0.05704157433714954
$ws.Range("T6").Value = 0.05734055747647147
$ws.Range("G7").Value = 10.311275
$ws.Range("H7").Value = 30.933825
$ws.Range("I7").Value = 0.5230601433626448
$ws.Range("J7").Value = 0.5230601433626447
$ws.Range("M7").Value = 1.332201
$ws.Range("N7").Value = 3.996603
$ws.Range("O7").Value = 0.4977000102551081
$ws.Range("P7").Value = 0.5003087024806611
$ws.Range("Q7").Value = 13.736690866275
$ws.Range("R7").Value = 123.630217796475
$ws.Range("S7").Value = 0.2603270387156266
$ws.Range("T7").Value = 0.2616915416451133
$ws.Range("G8").Value = 10.311275
$ws.Range("H8").Value = 30.933825
$ws.Range("I8").Value = 0.5230601433626448
$ws.Range("J8").Value = 0.5230601433626447
$ws.Range("O8").Value = 0.3400184143635273
$ws.Range("P8").Value = 0.341800619257677
$ws.Range("Q8").Value = 9.384624775391666
$ws.Range("R8").Value = 84.461622978525
$ws.Range("S8").Value = 0.1778500805629258
$ws.Range("T8").Value = 0.1787822809103612
$ws.Range("G9").Value = 10.311275
$ws.Range("H9").Value = 30.933825
$ws.Range("I9").Value = 0.5230601433626448
$ws.Range("J9").Value = 0.5230601433626447
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.07237833333333334
$ws.Range("N9").Value = 0.217135
$ws.Range("O9").Value = 0.02703998664034003
$ws.Range("P9").Value = 0.02718171660110807
$ws.Range("Q9").Value = 0.7463128990416668
$ws.Range("R9").Value = 6.716816091375
$ws.Range("S9").Value = 0.01414353928862026
$ws.Range("T9").Value = 0.01421767258221837
$ws.Range("G10").Value = 10.311275
$ws.Range("H10").Value = 30.933825
$ws.Range("I10").Value = 0.5230601433626448
$ws.Range("J10").Value = 0.5230601433626447
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.0418705
$ws.Range("N10").Value = 0.083741
$ws.Range("O10").Value = 0.01564249559892726
$ws.Range("P10").Value = 0.01048299044324218
$ws.Range("Q10").Value = 0.4317382398875
$ws.Range("R10").Value = 2.590429439325
$ws.Range("S10").Value = 0.008181965990524434
$ws.Range("T10").Value = 0.005483234484111491
$ws.Range("G11").Value = 10.311275
$ws.Range("H11").Value = 30.933825
$ws.Range("I11").Value = 0.5230601433626448
$ws.Range("J11").Value = 0.5230601433626447
$ws.Range("M11").Value = 0.3201326666666667
$ws.Range("N11").Value = 0.960398
$ws.Range("O11").Value = 0.1195990931420972
$ws.Range("P11").Value = 0.1202259712173118
$ws.Range("Q11").Value = 3.300975962483334
$ws.Range("R11").Value = 29.70878366235
$ws.Range("S11").Value = 0.06255751880494768
$ws.Range("T11").Value = 0.06288541374084029